# MIPS_By_Specialty_Summary.pptx edit script
# - Slide 1: new two-part title ("Top 15 Specialties by MIPS Score
#   Distribution " + red "2022 CMS Data"), disable title autofit, move the
#   chart picture down a bit.
# - Slide 2 (new, Title+Content layout): title mirrors slide 1 (bigger),
#   content placeholder holds three bullet paragraphs describing the chart.

$p = $ppt.ActivePresentation

# EMU <-> point helper. PowerPoint COM positions/sizes are expressed in
# points (1 pt = 12700 EMU); nudge by a hair so float error doesn't make
# the internal EMU conversion round down to one-less-than-target.
function EmuToPt([double]$emu) {
    return ($emu / 12700.0) + 0.00001
}

# ---------------------------------------------------------------------
# Slide 1
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)

$title1 = $s1.Shapes.Item(1)
$title1.TextFrame.AutoSize = 0   # ppAutoSizeNone -> <a:noAutofit/>

$tr1 = $title1.TextFrame.TextRange
$tr1.Text = "Top 15 Specialties by MIPS Score Distribution 2022 CMS Data"

$run1a = $tr1.Characters(1, 46)
$run1a.LanguageID = "en-GB"
$run1a.Font.Size = 32

$run1b = $tr1.Characters(47, 13)
$run1b.LanguageID = "en-GB"
$run1b.Font.Size = 24
$run1b.Font.Color.RGB = 255

$pic1 = $s1.Shapes.Item(2)
$pic1.Left = EmuToPt 457200
$pic1.Top = EmuToPt 1575585
$pic1.Width = EmuToPt 8001000
$pic1.Height = EmuToPt 4572000

# ---------------------------------------------------------------------
# Slide 2 (new slide, "Title and Content" layout)
# ---------------------------------------------------------------------
$s2 = $p.Slides.Add(2, 2)

$title2 = $s2.Shapes.Item(1)
$title2.TextFrame.AutoSize = 2   # ppAutoSizeTextToFitShape -> <a:normAutofit/>

$tr2 = $title2.TextFrame.TextRange
$tr2.Text = "Top 15 Specialties by MIPS Score Distribution 2022 CMS Data"

$run2a = $tr2.Characters(1, 46)
$run2a.LanguageID = "en-GB"
$run2a.Font.Size = 44

$run2b = $tr2.Characters(47, 13)
$run2b.LanguageID = "en-GB"
$run2b.Font.Size = 36
$run2b.Font.Color.RGB = 255

$body = $s2.Shapes.Item(2)
$body.Left = EmuToPt 457200
$body.Top = EmuToPt 2166425
$body.Width = EmuToPt 8229600
$body.Height = EmuToPt 3959738
$body.TextFrame.AutoSize = 2   # ppAutoSizeTextToFitShape -> <a:normAutofit/>

$trBody = $body.TextFrame.TextRange

# --- paragraph 1 ---
$trBody.Text = "This boxplot shows MIPS score distributions for the 15 most common specialties."
$trBody.LanguageID = "en-GB"

$p1r2 = $trBody.Characters(20, 25)
$p1r2.LanguageID = "en-GB"
$p1r2.Font.Bold = $true

$p1r4 = $trBody.Characters(53, 15)
$p1r4.LanguageID = "en-GB"
$p1r4.Font.Color.RGB = 255

# --- paragraph 2 ---
[void]$trBody.InsertAfter("`rBoxes represent IQRs, lines show medians, and dots indicate outliers.")
$trAll2 = $body.TextFrame.TextRange
$para2 = $trAll2.Paragraphs(2, 1)
$para2.LanguageID = "en-GB"

# --- paragraph 3 ---
[void]$trAll2.InsertAfter("`rSpecialties like General Surgery and Cardiology trend higher, while Emergency Medicine shows wider variability.")
$trAll3 = $body.TextFrame.TextRange
$para3 = $trAll3.Paragraphs(3, 1)
$para3.LanguageID = "en-GB"

$p3r2 = $para3.Characters(18, 16)
$p3r2.Font.Bold = $true
$p3r4 = $para3.Characters(38, 10)
$p3r4.Font.Bold = $true
$p3r6 = $para3.Characters(49, 12)
$p3r6.Font.Italic = $true
$p3r6.Font.Underline = $true
$p3r8 = $para3.Characters(69, 18)
$p3r8.Font.Bold = $true
$p3r10 = $para3.Characters(94, 17)
$p3r10.Font.Italic = $true
$p3r10.Font.Underline = $true
